$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.075.39"
$ws.Range("E2").Value = "  +0.85%  "

$ws.Range("D3").Value = "1.747.74"
$ws.Range("E3").Value = "  +0.41%  "

$ws.Range("D4").Value = "'1.001"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.00%  "

$ws.Range("D5").Value = "'234.88"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +4.95%  "

$ws.Range("D6").Value = "'1.000"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.03%  "

$ws.Range("D7").Value = "'0.5288"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +2.99%  "

$ws.Range("D8").Value = "'0.2800"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.50%  "

$ws.Range("D9").Value = "'0.06184"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.75%  "

$ws.Range("D10").Value = "1.745.64"
$ws.Range("E10").Value = "  -0.17%  "

$ws.Range("D11").Value = "'0.07169"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +3.25%  "

$ws.Range("D12").Value = "'15.41"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.79%  "

$ws.Range("D13").Value = "'0.6437"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +2.14%  "

$ws.Range("D14").Value = "'4.608"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.94%  "

$ws.Range("D15").Value = "'78.35"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +2.99%  "

$ws.Range("D16").Value = "'1.000"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.11%  "

$ws.Range("D17").Value = "'1.000"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.04%  "

$ws.Range("D18").Value = "25.998.11"
$ws.Range("E18").Value = "  +0.54%  "

$ws.Range("D19").Value = "'11.68"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +2.69%  "

$ws.Range("D20").Value = "'0.000006730"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.59%  "

$ws.Range("D21").Value = "1.967.44"
$ws.Range("E21").Value = "  +0.26%  "

$ws.Range("D22").Value = "'4.320"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +6.29%  "

$ws.Range("D23").Value = "'8.726"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +4.42%  "

$ws.Range("D24").Value = "'5.235"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +2.89%  "

$ws.Range("D25").Value = "'138.59"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.28%  "

$ws.Range("E26").Value = "  +1.06%  "

$ws.Range("D27").Value = "'15.29"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +2.43%  "

$ws.Range("D28").Value = "'1.803"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.32%  "

$ws.Range("D29").Value = "'104.37"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +2.22%  "

$ws.Range("D30").Value = "'0.08277"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.30%  "

$ws.Range("D31").Value = "'3.802"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +5.62%  "

$ws.Range("D32").Value = "'3.644"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +7.79%  "

$ws.Range("D33").Value = "'0.04575"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +4.97%  "

$ws.Range("D34").Value = "'2.643"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.73%  "

$ws.Range("D35").Value = "'1.009"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +4.76%  "

$ws.Range("D36").Value = "'0.6337"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +5.68%  "

$ws.Range("D37").Value = "'2.698"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.52%  "

$ws.Range("D38").Value = "'0.01607"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +4.08%  "

$ws.Range("D39").Value = "'1.962"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +3.90%  "

$ws.Range("D40").Value = "'0.9997"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.04%  "

$ws.Range("D41").Value = "'100.69"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.60%  "

$ws.Range("D42").Value = "'0.3928"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +3.36%  "

$ws.Range("D43").Value = "'0.7438"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +3.49%  "

$ws.Range("D44").Value = "'5.023"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +2.73%  "

$ws.Range("D45").Value = "'0.1144"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +4.73%  "

$ws.Range("D46").Value = "'6.352"
$ws.Range("D46").Style = "Normal"

$ws.Range("D47").Value = "'0.05345"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.83%  "

$ws.Range("D48").Value = "'30.93"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +4.85%  "

$ws.Range("D49").Value = "'54.29"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +4.48%  "

$ws.Range("D50").Value = "'7.628"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +2.77%  "

$ws.Range("D51").Value = "'0.3455"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.88%  "

